$d = $word.ActiveDocument

function Replace-ParagraphRuns {
    param(
        [string]$anchorText,
        [string]$runsXml
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Anchor text not found: " + $anchorText)
    }
    # Expand the found hit to cover its whole enclosing paragraph (wdParagraph = 4)
    $rng.Expand(4)

    $xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# --- First paragraph: "К недостаткам работы можно отнести ..." ---
$runs1 = '<w:r><w:t>К недостаткам работы</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> можно отнести</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>проведе</w:t></w:r>' + `
    '<w:r><w:t>ние исследования</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> с использованием только одного семейства компиляторов и средств интерпретации, однако данный недостаток не влияет на общее качество работы, а дает возможность для дальнейшего совершенствования средств эмуляции</w:t></w:r>' + `
    '<w:r><w:t>.</w:t></w:r>'

Replace-ParagraphRuns "недостаткам работы можно отнести" $runs1

# --- Second paragraph: "Тем не менее, как руководитель ВКРМ считаю, ..." ---
$runs2 = '<w:r><w:t>К</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">ак руководитель ВКРМ считаю, что представленная студентом работа </w:t></w:r>' + `
    '<w:r><w:t>выполнена в соответствии с</w:t></w:r>' + `
    '<w:r w:rsidR="00CB7EA6"><w:t xml:space="preserve"> требованиями к ВКРМ и</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> техническим заданием,</w:t></w:r>' + `
    '<w:r w:rsidR="00CB7EA6"><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r w:rsidR="006F4731"><w:t xml:space="preserve">заслуживает оценки </w:t></w:r>' + `
    '<w:r w:rsidR="008A069A"><w:t>«</w:t></w:r>' + `
    '<w:r w:rsidR="006F4731"><w:t>отлично</w:t></w:r>' + `
    '<w:r w:rsidR="008A069A"><w:t>»</w:t></w:r>' + `
    '<w:r w:rsidR="006F4731"><w:t xml:space="preserve">, а автор присвоения квалификации </w:t></w:r>' + `
    '<w:r w:rsidR="008A069A"><w:t>магистра</w:t></w:r>' + `
    '<w:r w:rsidR="006F4731"><w:t xml:space="preserve"> по направлению 11.0</w:t></w:r>' + `
    '<w:r w:rsidR="008A069A"><w:t>4</w:t></w:r>' + `
    '<w:r w:rsidR="006F4731"><w:t xml:space="preserve">.03 «Конструирование и технология электронных </w:t></w:r>' + `
    '<w:r w:rsidR="001B6B1F"><w:t>средств</w:t></w:r>' + `
    '<w:r w:rsidR="006F4731"><w:t>»</w:t></w:r>' + `
    '<w:r w:rsidR="007302F1"><w:t>.</w:t></w:r>' + `
    '<w:r w:rsidR="00CB7EA6"><w:t xml:space="preserve"> Студент Кутаев К.С. рекомендуется к зачислению в аспирантуру.</w:t></w:r>'

Replace-ParagraphRuns "Тем не менее, как руководитель ВКРМ считаю" $runs2
